$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new blank worksheet named "Sheet7" right after "OnBoardingStep7"
#    (i.e. before "LoginData"). This matches the tab order:
#    Signup, OnBoardingStep7, Sheet7, LoginData, Dashboard, EditYourProfile,
#    OnBoarding
# ---------------------------------------------------------------------------
$afterSheet = $wb.Worksheets.Item("OnBoardingStep7")
$newSheet = $wb.Worksheets.Add($null, $afterSheet)
$newSheet.Name = "Sheet7"

# ---------------------------------------------------------------------------
# 2. Add two new columns (validation hooks for login UI) to the
#    "OnBoardingStep7" sheet: "validDate" / "InvalidDate" headers in E1:F1,
#    with sample date values in E2:F2 formatted as dates.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("OnBoardingStep7")

$ws.Range("E1").Value = "validDate"
$ws.Range("F1").Value = "InvalidDate"

# Raw date serials (2025-02-14 and 2023-08-05) -- set numerically so no time
# fraction sneaks in.
$ws.Range("E2").Value = 45702.0
$ws.Range("F2").Value = 45143.0

# Visual styling: white fill + a small blue monospace font, each column
# carrying its own date format.
$dateRange = $ws.Range("E2:F2")
$dateRange.Interior.Color = 16777215
$dateFont = $dateRange.Font
$dateFont.Name = "Menlo"
$dateFont.Size = 12
$dateFont.Color = 16711722

$ws.Range("E2").NumberFormat = "m/d/yyyy"
$ws.Range("F2").NumberFormat = "mm/dd/yyyy"

# Restore the originally active sheet/selection so inserting the new sheet
# doesn't shift the workbook's focus.
$wb.Worksheets.Item("Signup").Activate() | Out-Null
$wb.Worksheets.Item("Signup").Range("A1").Select() | Out-Null
